$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value on a cell while forcing text interpretation
# (prevents Excel from auto-converting numeric-looking strings into numbers)
# and then restore the cell style so no stray number-format style is left behind.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '47.093.30'
Set-TextValue $ws.Range("E2") '  +1.12%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.485.50'
Set-TextValue $ws.Range("E3") '  +0.84%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '319.88'
Set-TextValue $ws.Range("E5") '  -1.05%  '

# Row 6
Set-TextValue $ws.Range("D6") '108.22'
Set-TextValue $ws.Range("E6") '  +2.99%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.521'
Set-TextValue $ws.Range("E7") '  -0.14%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.04%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.532'
Set-TextValue $ws.Range("E9") '  -1.00%  '

# Row 10
Set-TextValue $ws.Range("D10") '38.60'
Set-TextValue $ws.Range("E10") '  +7.06%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0808'
Set-TextValue $ws.Range("E11") '  -0.78%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +0.28%  '

# Row 13
Set-TextValue $ws.Range("D13") '18.15'
Set-TextValue $ws.Range("E13") '  -0.70%  '

# Row 14
Set-TextValue $ws.Range("D14") '7.10'
Set-TextValue $ws.Range("E14") '  +0.34%  '

# Row 15
Set-TextValue $ws.Range("D15") '2.877.53'
Set-TextValue $ws.Range("E15") '  +1.04%  '

# Row 16
Set-TextValue $ws.Range("D16") '2.489.21'
Set-TextValue $ws.Range("E16") '  -0.08%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.843'
Set-TextValue $ws.Range("E17") '  +0.01%  '

# Row 18
Set-TextValue $ws.Range("D18") '47.026.36'
Set-TextValue $ws.Range("E18") '  +1.23%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.67'
Set-TextValue $ws.Range("E19") '  +0.37%  '

# Row 20
Set-TextValue $ws.Range("D20") '6.58'
Set-TextValue $ws.Range("E20") '  +1.85%  '

# Row 21
Set-TextValue $ws.Range("D21") '2.75'
Set-TextValue $ws.Range("E21") '  +15.98%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.0₃0932'
Set-TextValue $ws.Range("E22") '  -0.37%  '

# Row 23
Set-TextValue $ws.Range("D23") '70.65'
Set-TextValue $ws.Range("E23") '  +0.32%  '

# Row 24
Set-TextValue $ws.Range("D24") '245.37'
Set-TextValue $ws.Range("E24") '  -1.35%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.55'
Set-TextValue $ws.Range("E25") '  +0.16%  '

# Row 26
Set-TextValue $ws.Range("E26") '  -0.01%  '

# Row 27
Set-TextValue $ws.Range("D27") '25.63'
Set-TextValue $ws.Range("E27") '  -1.75%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.28'
Set-TextValue $ws.Range("E28") '  -1.31%  '

# Row 29
Set-TextValue $ws.Range("D29") '10.01'
Set-TextValue $ws.Range("E29") '  +2.32%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +8.69%  '

# Row 31
Set-TextValue $ws.Range("D31") '34.84'
Set-TextValue $ws.Range("E31") '  -0.64%  '

# Row 32
Set-TextValue $ws.Range("D32") '49.91'
Set-TextValue $ws.Range("E32") '  +0.79%  '

# Row 33
Set-TextValue $ws.Range("D33") '20.00'
Set-TextValue $ws.Range("E33") '  +2.25%  '

# Row 34
Set-TextValue $ws.Range("D34") '5.33'
Set-TextValue $ws.Range("E34") '  -0.08%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0780'
Set-TextValue $ws.Range("E35") '  +1.67%  '

# Row 36
Set-TextValue $ws.Range("E36") '  +0.22%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.95'
Set-TextValue $ws.Range("E37") '  +2.52%  '

# Row 38
Set-TextValue $ws.Range("D38") '4.64'
Set-TextValue $ws.Range("E38") '  +0.27%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.94'
Set-TextValue $ws.Range("E39") '  +0.18%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +0.05%  '

# Row 41
Set-TextValue $ws.Range("E41") '  -0.40%  '

# Row 42
Set-TextValue $ws.Range("D42") '119.21'
Set-TextValue $ws.Range("E42") '  -3.34%  '

# Row 43
Set-TextValue $ws.Range("D43") '21.30'
Set-TextValue $ws.Range("E43") '  +2.17%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.0293'
Set-TextValue $ws.Range("E44") '  +0.01%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.976.86'
Set-TextValue $ws.Range("E45") '  -0.30%  '

# Row 46
Set-TextValue $ws.Range("D46") '3.00'
Set-TextValue $ws.Range("E46") '  +1.04%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.00'
Set-TextValue $ws.Range("E47") '  -3.31%  '

# Row 48
Set-TextValue $ws.Range("D48") '9.04'
Set-TextValue $ws.Range("E48") '  +0.96%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.75'
Set-TextValue $ws.Range("E49") '  -2.45%  '

# Row 50
Set-TextValue $ws.Range("E50") '  -4.07%  '

# Row 51
Set-TextValue $ws.Range("D51") '57.03'
Set-TextValue $ws.Range("E51") '  +4.86%  '
